$d = $word.ActiveDocument

# 1. Merge "Integrantes do Projeto – Turma SI" + "." + "B" into one run's text.
$d.Content.Find.Execute("Integrantes do Projeto – Turma SI.B", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Integrantes do Projeto – Turma SI.B", 2)

# 4. "Nome da Equipe de Desenvolvimento: " -> "Equipe de Desenvolvimento: "
$d.Content.Find.Execute("Nome da Equipe de Desenvolvimento: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Equipe de Desenvolvimento: ", 2)
